# Auto-generated cell value updates based on the provided OOXML diff.
# Applies updated market-price / profit figures to the ALC, ARM, BSM, CRP,
# CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 193.88889
$ws.Range("I11").Value = 193.88889
$ws.Range("K11").Value = 193.88889
$ws.Range("M11").Value = -53.88889
$ws.Range("H40").Value = 1286.15
$ws.Range("I40").Value = 772.3
$ws.Range("J40").Value = 1800
$ws.Range("K40").Value = 772.3
$ws.Range("L40").Value = 1800
$ws.Range("M40").Value = -597.3
$ws.Range("N40").Value = -2150
$ws.Range("H112").Value = 3788928.8
$ws.Range("J112").Value = 1101.5161
$ws.Range("L112").Value = 3304.5483
$ws.Range("N112").Value = -5520.5483
$ws.Range("H129").Value = 855.6901
$ws.Range("J129").Value = 869.597
$ws.Range("L129").Value = 2608.791
$ws.Range("N129").Value = -12608.791
$ws.Range("H132").Value = 59585.168
$ws.Range("I132").Value = 63075.766
$ws.Range("K132").Value = 189227.298
$ws.Range("M132").Value = -186697.298

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3269.0667
$ws.Range("I45").Value = 2952
$ws.Range("K45").Value = 2952
$ws.Range("M45").Value = -2575
$ws.Range("H61").Value = 2279.28
$ws.Range("I61").Value = 1004.9286
$ws.Range("J61").Value = 3901.182
$ws.Range("K61").Value = 1004.9286
$ws.Range("L61").Value = 3901.182
$ws.Range("M61").Value = -792.9286
$ws.Range("N61").Value = -4325.182
$ws.Range("H74").Value = 861.1852
$ws.Range("I74").Value = 439.57895
$ws.Range("J74").Value = 1862.5
$ws.Range("K74").Value = 439.57895
$ws.Range("L74").Value = 1862.5
$ws.Range("M74").Value = 434.42105
$ws.Range("N74").Value = -3610.5
$ws.Range("H77").Value = 861.1852
$ws.Range("I77").Value = 439.57895
$ws.Range("J77").Value = 1862.5
$ws.Range("K77").Value = 2197.89475
$ws.Range("L77").Value = 9312.5
$ws.Range("M77").Value = 2170.10525
$ws.Range("N77").Value = -18048.5
$ws.Range("H136").Value = 2279.28
$ws.Range("I136").Value = 1004.9286
$ws.Range("J136").Value = 3901.182
$ws.Range("K136").Value = 3014.7858
$ws.Range("L136").Value = 11703.546
$ws.Range("M136").Value = -464.7857999999997
$ws.Range("N136").Value = -16803.546

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 916.6667
$ws.Range("I20").Value = 880
$ws.Range("J20").Value = 1100
$ws.Range("K20").Value = 880
$ws.Range("L20").Value = 1100
$ws.Range("M20").Value = -633
$ws.Range("N20").Value = -1594
$ws.Range("H80").Value = 952.1539
$ws.Range("I80").Value = 1706
$ws.Range("J80").Value = 617.1111
$ws.Range("K80").Value = 1706
$ws.Range("L80").Value = 617.1111
$ws.Range("M80").Value = -708
$ws.Range("N80").Value = -2613.1111
$ws.Range("H83").Value = 952.1539
$ws.Range("I83").Value = 1706
$ws.Range("J83").Value = 617.1111
$ws.Range("K83").Value = 8530
$ws.Range("L83").Value = 3085.5555
$ws.Range("M83").Value = -3538
$ws.Range("N83").Value = -13069.5555

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17107.033
$ws.Range("I58").Value = 1049.96
$ws.Range("J58").Value = 84011.5
$ws.Range("K58").Value = 1049.96
$ws.Range("L58").Value = 84011.5
$ws.Range("M58").Value = -846.96
$ws.Range("N58").Value = -84417.5
$ws.Range("H94").Value = 5997.5884
$ws.Range("J94").Value = 7932.636
$ws.Range("L94").Value = 7932.636
$ws.Range("N94").Value = -8834.636
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 22893.809
$ws.Range("I132").Value = 28463.842
$ws.Range("J132").Value = 7775.143
$ws.Range("K132").Value = 85391.526
$ws.Range("L132").Value = 23325.429
$ws.Range("M132").Value = -82861.526
$ws.Range("N132").Value = -28385.429
$ws.Range("H134").Value = 576.7895
$ws.Range("I134").Value = 561.82355
$ws.Range("J134").Value = 704
$ws.Range("K134").Value = 1685.47065
$ws.Range("L134").Value = 2112
$ws.Range("M134").Value = 849.5293500000002
$ws.Range("N134").Value = -7182
$ws.Range("H136").Value = 17107.033
$ws.Range("I136").Value = 1049.96
$ws.Range("J136").Value = 84011.5
$ws.Range("K136").Value = 3149.88
$ws.Range("L136").Value = 252034.5
$ws.Range("M136").Value = -599.8800000000001
$ws.Range("N136").Value = -257134.5
$ws.Range("H137").Value = 51180
$ws.Range("J137").Value = 51180
$ws.Range("L137").Value = 51180
$ws.Range("N137").Value = -61380

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 675.7838
$ws.Range("I5").Value = 523.58826
$ws.Range("K5").Value = 1570.76478
$ws.Range("M5").Value = -1458.76478
$ws.Range("H68").Value = 1362.45
$ws.Range("J68").Value = 1438.8889
$ws.Range("L68").Value = 4316.6667
$ws.Range("N68").Value = -5938.6667
$ws.Range("H71").Value = 1362.45
$ws.Range("J71").Value = 1438.8889
$ws.Range("L71").Value = 12950.0001
$ws.Range("N71").Value = -21062.0001
$ws.Range("H92").Value = 1056.8572
$ws.Range("I92").Value = 650
$ws.Range("J92").Value = 1219.6
$ws.Range("K92").Value = 1950
$ws.Range("L92").Value = 3658.8
$ws.Range("M92").Value = -702
$ws.Range("N92").Value = -6154.799999999999
$ws.Range("H131").Value = 157094.95
$ws.Range("J131").Value = 170325.14
$ws.Range("L131").Value = 510975.42
$ws.Range("N131").Value = -521055.42
$ws.Range("H135").Value = 675.7838
$ws.Range("I135").Value = 523.58826
$ws.Range("K135").Value = 4712.29434
$ws.Range("M135").Value = -2177.29434
$ws.Range("H140").Value = 4859.759
$ws.Range("I140").Value = 5849.421
$ws.Range("J140").Value = 2979.4
$ws.Range("K140").Value = 17548.263
$ws.Range("L140").Value = 8938.200000000001
$ws.Range("M140").Value = -12368.263
$ws.Range("N140").Value = -19298.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2150.5217
$ws.Range("I102").Value = 2140.8
$ws.Range("K102").Value = 2140.8
$ws.Range("M102").Value = -518.8000000000002
$ws.Range("H126").Value = 4391.4287
$ws.Range("I126").Value = 3392.2222
$ws.Range("K126").Value = 10176.6666
$ws.Range("M126").Value = -7706.6666
$ws.Range("H132").Value = 62220.72
$ws.Range("I132").Value = 51776
$ws.Range("J132").Value = 103999.6
$ws.Range("K132").Value = 155328
$ws.Range("L132").Value = 311998.8
$ws.Range("M132").Value = -152798
$ws.Range("N132").Value = -317058.8

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1027.238
$ws.Range("I16").Value = 804.7059
$ws.Range("K16").Value = 804.7059
$ws.Range("M16").Value = -634.7059
$ws.Range("H93").Value = 2293.389
$ws.Range("I93").Value = 2058.4666
$ws.Range("K93").Value = 2058.4666
$ws.Range("M93").Value = -810.4666000000002
$ws.Range("H136").Value = 23911.39
$ws.Range("I136").Value = 32297.625
$ws.Range("J136").Value = 4742.857
$ws.Range("K136").Value = 96892.875
$ws.Range("L136").Value = 14228.571
$ws.Range("M136").Value = -94342.875
$ws.Range("N136").Value = -19328.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1898782.1
$ws.Range("I136").Value = 2688933.2
$ws.Range("K136").Value = 8066799.600000001
$ws.Range("M136").Value = -8064249.600000001
